$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '62.238.04'
$cell.Style = "Normal"
$ws.Range("E2").Value = '  -0.33%  '

# Row 3
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '2.468.47'
$cell.Style = "Normal"
$ws.Range("E3").Value = '  +1.69%  '

# Row 4
$ws.Range("E4").Value = '  -0.15%  '

# Row 5
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '584.21'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  +1.16%  '

# Row 6
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '142.76'
$cell.Style = "Normal"
$ws.Range("E6").Value = '  -0.21%  '

# Row 7
$ws.Range("E7").Value = '  -0.08%  '

# Row 8
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '0.531'
$cell.Style = "Normal"
$ws.Range("E8").Value = '  +0.95%  '

# Row 9
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '2.463.69'
$cell.Style = "Normal"
$ws.Range("E9").Value = '  +1.44%  '

# Row 10
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '0.111'
$cell.Style = "Normal"
$ws.Range("E10").Value = '  +4.37%  '

# Row 11
$ws.Range("E11").Value = '  +3.05%  '

# Row 12
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '5.20'
$cell.Style = "Normal"
$ws.Range("E12").Value = '  +0.17%  '

# Row 13
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '0.343'
$cell.Style = "Normal"
$ws.Range("E13").Value = '  -0.94%  '

# Row 14
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '26.14'
$cell.Style = "Normal"
$ws.Range("E14").Value = '  -0.98%  '

# Row 15
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '0.0000174'
$cell.Style = "Normal"
$ws.Range("E15").Value = '  +0.74%  '

# Row 16
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '2.897.80'
$cell.Style = "Normal"
$ws.Range("E16").Value = '  +0.56%  '

# Row 17
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '62.126.19'
$cell.Style = "Normal"
$ws.Range("E17").Value = '  -0.52%  '

# Row 18
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '2.464.10'
$cell.Style = "Normal"
$ws.Range("E18").Value = '  +1.70%  '

# Row 19
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '10.72'
$cell.Style = "Normal"
$ws.Range("E19").Value = '  -2.26%  '

# Row 20
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '7.37'
$cell.Style = "Normal"
$ws.Range("E20").Value = '  +4.00%  '

# Row 21
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '327.43'
$cell.Style = "Normal"
$ws.Range("E21").Value = '  -0.97%  '

# Row 22
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '4.10'
$cell.Style = "Normal"
$ws.Range("E22").Value = '  -0.29%  '

# Row 23
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '1.95'
$cell.Style = "Normal"
$ws.Range("E23").Value = '  -1.40%  '

# Row 24
$ws.Range("E24").Value = '  -0.07%  '

# Row 25
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '65.41'
$cell.Style = "Normal"
$ws.Range("E25").Value = '  -0.61%  '

# Row 26
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '9.27'
$cell.Style = "Normal"
$ws.Range("E26").Value = '  +3.48%  '

# Row 27
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '588.18'
$cell.Style = "Normal"
$ws.Range("E27").Value = '  -6.78%  '

# Row 28
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '2.569.33'
$cell.Style = "Normal"

# Row 29
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '0.998'
$cell.Style = "Normal"
$ws.Range("E29").Value = '  -0.24%  '

# Row 30
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '0.0₃0945'
$cell.Style = "Normal"
$ws.Range("E30").Value = '  -0.78%  '

# Row 31
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '8.00'
$cell.Style = "Normal"
$ws.Range("E31").Value = '  +0.02%  '

# Row 32
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '1.39'
$cell.Style = "Normal"
$ws.Range("E32").Value = '  -2.21%  '

# Row 33
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '1.89'
$cell.Style = "Normal"
$ws.Range("E33").Value = '  +0.64%  '

# Row 34
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '0.135'
$cell.Style = "Normal"
$ws.Range("E34").Value = '  -3.40%  '

# Row 35
$ws.Range("E35").Value = '  -0.22%  '

# Row 36
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '4.81'
$cell.Style = "Normal"
$ws.Range("E36").Value = '  -3.28%  '

# Row 37
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '1.42'
$cell.Style = "Normal"
$ws.Range("E37").Value = '  -1.04%  '

# Row 38
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '0.374'
$cell.Style = "Normal"
$ws.Range("E38").Value = '  -0.14%  '

# Row 39
$ws.Range("B39").Value = 'Monero'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '152.44'
$cell.Style = "Normal"
$ws.Range("E39").Value = '  +2.24%  '

# Row 40
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '18.41'
$cell.Style = "Normal"
$ws.Range("E40").Value = '  +0.14%  '

# Row 41
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '5.23'
$cell.Style = "Normal"
$ws.Range("E41").Value = '  -0.38%  '

# Row 42
$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '1.71'
$cell.Style = "Normal"
$ws.Range("E42").Value = '  -1.77%  '

# Row 43
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '42.26'
$cell.Style = "Normal"
$ws.Range("E43").Value = '  -0.40%  '

# Row 44
$ws.Range("B44").Value = 'USDe'
$ws.Range("C44").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '0.999'
$cell.Style = "Normal"
$ws.Range("E44").Value = '  +0.04%  '

# Row 45
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '2.41'
$cell.Style = "Normal"
$ws.Range("E45").Value = '  -1.81%  '

# Row 46
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '0.0₆0294'
$cell.Style = "Normal"
$ws.Range("E46").Value = '  +23.19%  '

# Row 47
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '143.64'
$cell.Style = "Normal"
$ws.Range("E47").Value = '  +0.20%  '

# Row 48
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '3.60'
$cell.Style = "Normal"
$ws.Range("E48").Value = '  -1.94%  '

# Row 49
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '0.604'
$cell.Style = "Normal"
$ws.Range("E49").Value = '  +1.37%  '

# Row 50
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '20.00'
$cell.Style = "Normal"
$ws.Range("E50").Value = '  +3.16%  '

# Row 51
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '0.0515'
$cell.Style = "Normal"
$ws.Range("E51").Value = '  -0.81%  '

Write-Host "Update complete"